# Add a new "video link" column (F) to the manage_product sheet:
#   F1 header "video link" (same header styling as the rest of row 1)
#   F2 a shutterstock video URL, word-wrapped like the other description cell
#   Column F sized to fit the long URL.
#
# The three other sheets (subcategory, footertext) are untouched by this
# script; their shared-string index shift in the saved file is a natural
# side effect of the two new strings being added to the workbook's shared
# string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manage_product")

# --- Header cell F1: same text + formatting as the rest of the header row ---
$ws.Range("F1").Value = "video link"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats - keep the value we just set

# --- Data cell F2: the video URL, wrapped like the neighbouring description cell ---
$ws.Range("F2").Value = "https://www.shutterstock.com/video/clip-3557755893-retro-projector-style-5-seconds-countdown-ten"
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats - base formatting (same font as the row)
$ws.Range("F2").WrapText = $true

# --- Column F width, sized for the URL text ---
$ws.Columns.Item(6).ColumnWidth = 25.33
